$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$r = $ws.Range("Z2")
$r.Borders.LineStyle = 1
# try leaving color completely untouched (default) --- that's what we currently do, confirmed no <color> emitted
